$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = 20.03
$ws.Cells.Item(2, 3).Value = 21.67
$ws.Cells.Item(2, 4).Value = 19.34
$ws.Cells.Item(2, 5).Value = 18.66

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 3).Value = 2.63
$ws.Cells.Item(2, 4).Value = 7.89
$ws.Cells.Item(2, 5).Value = 13.16
$ws.Cells.Item(2, 6).Value = 76.31999999999999
$ws.Cells.Item(3, 1).Value = 'Jharkhand'
$ws.Cells.Item(3, 3).Value = 8.33
$ws.Cells.Item(3, 4).Value = 8.33
$ws.Cells.Item(3, 5).Value = 29.17
$ws.Cells.Item(3, 6).Value = 54.17
$ws.Cells.Item(3, 7).Value = 'IN-JH'
$ws.Cells.Item(4, 1).Value = 'Puducherry'
$ws.Cells.Item(4, 3).Value = 0
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 25
$ws.Cells.Item(4, 6).Value = 50
$ws.Cells.Item(4, 7).Value = 'IN-PY'
$ws.Cells.Item(5, 4).Value = 50
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 50
$ws.Cells.Item(6, 1).Value = 'Mizoram'
$ws.Cells.Item(6, 3).Value = 18.18
$ws.Cells.Item(6, 4).Value = 18.18
$ws.Cells.Item(6, 5).Value = 18.18
$ws.Cells.Item(6, 6).Value = 45.45
$ws.Cells.Item(6, 7).Value = 'IN-MZ'
$ws.Cells.Item(7, 3).Value = 8
$ws.Cells.Item(7, 4).Value = 18.67
$ws.Cells.Item(7, 5).Value = 28
$ws.Cells.Item(7, 6).Value = 45.33
$ws.Cells.Item(8, 3).Value = 12
$ws.Cells.Item(8, 4).Value = 20
$ws.Cells.Item(8, 5).Value = 8
$ws.Cells.Item(8, 6).Value = 40
$ws.Cells.Item(9, 1).Value = 'Nagaland'
$ws.Cells.Item(9, 3).Value = 9.09
$ws.Cells.Item(9, 4).Value = 36.36
$ws.Cells.Item(9, 5).Value = 9.09
$ws.Cells.Item(9, 6).Value = 36.36
$ws.Cells.Item(9, 7).Value = 'IN-NL'
$ws.Cells.Item(10, 1).Value = 'Manipur'
$ws.Cells.Item(10, 3).Value = 12.5
$ws.Cells.Item(10, 4).Value = 25
$ws.Cells.Item(10, 5).Value = 31.25
$ws.Cells.Item(10, 6).Value = 31.25
$ws.Cells.Item(10, 7).Value = 'IN-MN'
$ws.Cells.Item(11, 3).Value = 9.09
$ws.Cells.Item(11, 4).Value = 22.73
$ws.Cells.Item(11, 5).Value = 40.91
$ws.Cells.Item(11, 6).Value = 27.27
$ws.Cells.Item(12, 1).Value = 'Odisha'
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 16.67
$ws.Cells.Item(12, 4).Value = 26.67
$ws.Cells.Item(12, 5).Value = 33.33
$ws.Cells.Item(12, 6).Value = 16.67
$ws.Cells.Item(12, 7).Value = 'IN-OR'
$ws.Cells.Item(13, 1).Value = 'Telangana'
$ws.Cells.Item(13, 3).Value = 24.24
$ws.Cells.Item(13, 4).Value = 33.33
$ws.Cells.Item(13, 5).Value = 18.18
$ws.Cells.Item(13, 6).Value = 15.15
$ws.Cells.Item(13, 7).Value = 'IN-TS'
$ws.Cells.Item(14, 1).Value = 'Tripura'
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 25
$ws.Cells.Item(14, 5).Value = 62.5
$ws.Cells.Item(14, 6).Value = 12.5
$ws.Cells.Item(14, 7).Value = 'IN-TR'
$ws.Cells.Item(15, 1).Value = 'Meghalaya'
$ws.Cells.Item(15, 3).Value = 36.36
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 54.55
$ws.Cells.Item(15, 6).Value = 9.09
$ws.Cells.Item(15, 7).Value = 'IN-ML'
$ws.Cells.Item(16, 1).Value = 'Delhi'
$ws.Cells.Item(16, 2).Value = 0
$ws.Cells.Item(16, 3).Value = 18.18
$ws.Cells.Item(16, 4).Value = 27.27
$ws.Cells.Item(16, 5).Value = 45.45
$ws.Cells.Item(16, 6).Value = 9.09
$ws.Cells.Item(16, 7).Value = 'IN-DL'
$ws.Cells.Item(17, 1).Value = 'Jammu and Kashmir'
$ws.Cells.Item(17, 3).Value = 9.09
$ws.Cells.Item(17, 4).Value = 40.91
$ws.Cells.Item(17, 5).Value = 40.91
$ws.Cells.Item(17, 6).Value = 9.09
$ws.Cells.Item(17, 7).Value = 'IN-JK'
$ws.Cells.Item(18, 1).Value = 'Haryana'
$ws.Cells.Item(18, 3).Value = 13.64
$ws.Cells.Item(18, 4).Value = 40.91
$ws.Cells.Item(18, 5).Value = 31.82
$ws.Cells.Item(18, 6).Value = 9.09
$ws.Cells.Item(18, 7).Value = 'IN-HR'
$ws.Cells.Item(19, 1).Value = 'Assam'
$ws.Cells.Item(19, 3).Value = 15.15
$ws.Cells.Item(19, 4).Value = 30.3
$ws.Cells.Item(19, 5).Value = 18.18
$ws.Cells.Item(19, 6).Value = 9.09
$ws.Cells.Item(19, 7).Value = 'IN-AS'
$ws.Cells.Item(20, 3).Value = 11.54
$ws.Cells.Item(20, 4).Value = 40.38
$ws.Cells.Item(20, 5).Value = 30.77
$ws.Cells.Item(20, 6).Value = 7.69
$ws.Cells.Item(21, 3).Value = 34.78
$ws.Cells.Item(21, 4).Value = 26.09
$ws.Cells.Item(21, 5).Value = 8.699999999999999
$ws.Cells.Item(21, 6).Value = 4.35
$ws.Cells.Item(22, 3).Value = 25.93
$ws.Cells.Item(22, 4).Value = 40.74
$ws.Cells.Item(22, 5).Value = 11.11
$ws.Cells.Item(22, 6).Value = 3.7
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 30.3
$ws.Cells.Item(23, 4).Value = 27.27
$ws.Cells.Item(23, 5).Value = 18.18
$ws.Cells.Item(23, 6).Value = 3.03
$ws.Cells.Item(24, 3).Value = 41.67
$ws.Cells.Item(24, 4).Value = 16.67
$ws.Cells.Item(24, 5).Value = 16.67
$ws.Cells.Item(25, 3).Value = 54.55
$ws.Cells.Item(25, 4).Value = 9.09
$ws.Cells.Item(25, 5).Value = 6.06
$ws.Cells.Item(26, 3).Value = 30.56
$ws.Cells.Item(26, 4).Value = 19.44
$ws.Cells.Item(26, 5).Value = 5.56
$ws.Cells.Item(27, 3).Value = 26.67
$ws.Cells.Item(27, 4).Value = 6.67
$ws.Cells.Item(27, 5).Value = 3.33
$ws.Cells.Item(28, 2).Value = 0
$ws.Cells.Item(28, 3).Value = 40.54
$ws.Cells.Item(28, 4).Value = 5.41
$ws.Cells.Item(28, 5).Value = 2.7
$ws.Cells.Item(29, 4).Value = 100
$ws.Cells.Item(30, 4).Value = 100
$ws.Cells.Item(31, 2).Value = 0
$ws.Cells.Item(31, 3).Value = 53.85
$ws.Cells.Item(31, 4).Value = 15.38
$ws.Cells.Item(32, 2).Value = 0
$ws.Cells.Item(32, 3).Value = 23.08
$ws.Cells.Item(32, 4).Value = 7.69
$ws.Cells.Item(33, 3).Value = 50

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 3).Value = 2.22
$ws.Cells.Item(2, 4).Value = 6.67
$ws.Cells.Item(2, 5).Value = 15.56
$ws.Cells.Item(2, 6).Value = 75.56
$ws.Cells.Item(3, 2).Value = 0
$ws.Cells.Item(3, 3).Value = 14.78
$ws.Cells.Item(3, 4).Value = 23.48
$ws.Cells.Item(3, 5).Value = 23.48
$ws.Cells.Item(3, 6).Value = 25.22
